$d = $word.ActiveDocument
$sec = $d.Sections.First

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
# Footers.Item(1) -> footer2.xml ("default" footer), Footers.Item(2) -> footer1.xml ("first page" footer)
# Headers.Item(1) -> header2.xml ("default" header), Headers.Item(2) -> header1.xml ("first page" header)

# footer1.xml / footer2.xml: Pearson logo, rename image2.png -> image1.png
$f1 = $sec.Footers.Item(2)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
    $shp = $f1.Range.InlineShapes.Item(1)
    $shp.Name = "image1.png"
}

$f2 = $sec.Footers.Item(1)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
    $shp = $f2.Range.InlineShapes.Item(1)
    $shp.Name = "image1.png"
}

# header1.xml / header2.xml: BTEC logo, rename image1.jpg -> image2.jpg
$h1 = $sec.Headers.Item(2)
if ($h1.Exists -and $h1.Range.InlineShapes.Count -ge 1) {
    $shp = $h1.Range.InlineShapes.Item(1)
    $shp.Name = "image2.jpg"
}

$h2 = $sec.Headers.Item(1)
if ($h2.Exists -and $h2.Range.InlineShapes.Count -ge 1) {
    $shp = $h2.Range.InlineShapes.Item(1)
    $shp.Name = "image2.jpg"
}
